# "Prepare for work without queue"
# - Remove the Orchestrator-queue related rows from the Settings sheet
#   (OrchestratorQueueName / OrchestratorQueueFolder), which shifts the
#   logF_BusinessProcessName row up.
# - Rename the "Framework" business-process value to the project's own
#   name, Unicorn_Name_Generator.
# - Update a couple of Constants values (MaxRetryNumber,
#   MaxConsecutiveSystemExceptions, ShouldMarkJobAsFaulted).
# - Switch the active tab from Assets to Constants.

$wb = $excel.ActiveWorkbook

# --- Settings sheet: drop the queue-name / queue-folder rows -------------
$ws1 = $wb.Worksheets.Item("Settings")

$ws1.Rows("2:2").Select() | Out-Null
$ws1.Rows("2:2").Delete() | Out-Null
$ws1.Rows("2:2").Delete() | Out-Null
$ws1.Rows("2:2").Delete() | Out-Null

# logF_BusinessProcessName (now row 2) keeps its Name/Description, only the
# Value column changes.
$ws1.Range("B2").Value = "Unicorn_Name_Generator"

# --- Constants sheet: update default values, then make it the active tab -
$ws2 = $wb.Worksheets.Item("Constants")

$ws2.Range("B2").Value = 3
$ws2.Range("B3").Value = 4
$ws2.Range("B17").Value = $true

$ws2.Activate() | Out-Null
$ws2.Range("B20").Select() | Out-Null
